$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: move footer block (___ signature lines) from rows 26-27 down to rows 43-44 ---
$ws.Range("B26:C27").UnMerge()
$ws.Range("H26:J27").UnMerge()

$ws.Range("B26:J27").Copy()
$ws.Range("B43:J44").PasteSpecial(-4122)
$ws.Range("B43:J44").PasteSpecial(-4163)
$excel.CutCopyMode = 0

$ws.Range("B43:C43").Merge()
$ws.Range("B44:C44").Merge()
$ws.Range("H43:J43").Merge()
$ws.Range("H44:J44").Merge()

# --- Step 2: extend the data-table formatting down to hold 23 worker rows (16-38) ---
# Row 21 (old) carries the thicker bottom border used for the LAST row of the table;
# grab that pattern first and stash it on new row 38 before it gets overwritten.
$ws.Range("B21:J21").Copy()
$ws.Range("B38:J38").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 20 (old) carries the plain interior-row formatting; stamp it across rows 21-37.
$ws.Range("B20:J20").Copy()
$ws.Range("B21:J37").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Step 3: write the 23 worker rows of data ---
$data = @(
    @("CC","3800675","JAIRO JESUS SANCHEZ MENDOZA","2308",46400,1160000),
    @("CC","3800675","JAIRO JESUS SANCHEZ MENDOZA","2307",46400,1160000),
    @("CC","3800675","JAIRO JESUS SANCHEZ MENDOZA","2306",46400,1160000),
    @("CC","1050968367","EUDES DE JESUS GUTIERREZ CERDA","2308",21654,1000000),
    @("CC","1043963104","FRANKLIN BALCEIRO PEREZ","2307",36341,1160000),
    @("CC","1043963104","FRANKLIN BALCEIRO PEREZ","2308",46400,1160000),
    @("CC","1043963104","FRANKLIN BALCEIRO PEREZ","2307",46400,1160000),
    @("CC","1043963104","FRANKLIN BALCEIRO PEREZ","2304",46400,1160000),
    @("CC","1043963104","FRANKLIN BALCEIRO PEREZ","2303",46400,1160000),
    @("CC","1043963104","FRANKLIN BALCEIRO PEREZ","2302",46400,1160000),
    @("CC","9186163","ARNOLDO ARZUZA PEDRAZA","2307",46400,1000000),
    @("CC","9186163","ARNOLDO ARZUZA PEDRAZA","2306",13920,1000000),
    @("CC","1143380904","MANUEL SALVADOR VILLA CARO","2307",21654,908526),
    @("CC","1047451793","ROIMA JULIO BLANCO","2308",46400,1160000),
    @("CC","1047451793","ROIMA JULIO BLANCO","2307",46400,1160000),
    @("CC","1047451793","ROIMA JULIO BLANCO","2306",46400,1160000),
    @("CC","1047451793","ROIMA JULIO BLANCO","2304",46400,1160000),
    @("CC","1047451793","ROIMA JULIO BLANCO","2303",46400,1160000),
    @("CC","1047451793","ROIMA JULIO BLANCO","2302",46400,1160000),
    @("CC","1050967361","JUAN CAMILO HERNANDEZ VIGGIANI","2307",21654,1000000),
    @("CC","1007959494","JAIRO ENRIQUE DIAZ FIORILLO","2308",46400,1160000),
    @("CC","1007959494","JAIRO ENRIQUE DIAZ FIORILLO","2307",46400,1160000),
    @("CC","1007959494","JAIRO ENRIQUE DIAZ FIORILLO","2306",13920,1160000)
)

$r = 16
foreach ($row in $data) {
    $ws.Range("B$r").Value = $row[0]
    $ws.Range("C$r").Value = $row[1]
    $ws.Range("D$r").Value = $row[2]
    $ws.Range("E$r").Value = $row[3]
    $ws.Range("F$r").Value = $row[4]
    $ws.Range("G$r").Value = $row[5]
    $r = $r + 1
}

# --- Step 4: header / summary cells ---
$ws.Range("E11").Value = 917943      # VALOR MORA total
$ws.Range("C13").Value = 8           # Cant. Trabajadores
$ws.Range("F13").Value = 6           # Cant. Periodos

# --- Step 5: widen column D so the longer worker names fit (bestFit) ---
$ws.Columns("D").ColumnWidth = 33.90625

Write-Host "edit complete"
